# Apply odds updates to "Jogos da Semana" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9
$ws.Range("G9").Value = 1.7
$ws.Range("H9").Value = 3.4
$ws.Range("J9").Value = 1.08
$ws.Range("K9").Value = 8
$ws.Range("L9").Value = 1.4
$ws.Range("M9").Value = 2.75
$ws.Range("N9").Value = 2.25
$ws.Range("O9").Value = 1.62
$ws.Range("P9").Value = 1.44
$ws.Range("Q9").Value = 2.63
$ws.Range("R9").Value = 2.2
$ws.Range("T9").Value = 5.5
$ws.Range("U9").Value = 7
$ws.Range("V9").Value = 9
$ws.Range("W9").Value = 12
$ws.Range("X9").Value = 17
$ws.Range("Y9").Value = 34
$ws.Range("Z9").Value = 7.5
$ws.Range("AA9").Value = 7
$ws.Range("AB9").Value = 21
$ws.Range("AC9").Value = 67
$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 26
$ws.Range("AF9").Value = 19
$ws.Range("AG9").Value = 51
$ws.Range("AH9").Value = 41
$ws.Range("AI9").Value = 51

# Row 10
$ws.Range("G10").Value = 1.45
$ws.Range("I10").Value = 6
$ws.Range("L10").Value = 1.18
$ws.Range("M10").Value = 4.5
$ws.Range("N10").Value = 1.6
$ws.Range("O10").Value = 2.3
$ws.Range("R10").Value = 1.73
$ws.Range("S10").Value = 2
$ws.Range("U10").Value = 8
$ws.Range("Z10").Value = 15
$ws.Range("AA10").Value = 9

# Row 14
$ws.Range("J14").Value = 1.11
$ws.Range("K14").Value = 6.5
$ws.Range("Z14").Value = 6.5

# Row 17
$ws.Range("J17").Value = 1.05
$ws.Range("K17").Value = 11
$ws.Range("L17").Value = 1.29
$ws.Range("M17").Value = 3.5
$ws.Range("N17").Value = 1.9
$ws.Range("O17").Value = 1.9
$ws.Range("P17").Value = 1.36
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 1.7
$ws.Range("S17").Value = 2.05
$ws.Range("T17").Value = 9.5
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 11
$ws.Range("AD17").Value = 9
$ws.Range("AE17").Value = 12
$ws.Range("AJ17").Value = 201

# Row 18
$ws.Range("G18").Value = 2.8
$ws.Range("H18").Value = 3.25
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11
$ws.Range("L18").Value = 1.25
$ws.Range("M18").Value = 3.75
$ws.Range("N18").Value = 1.9
$ws.Range("O18").Value = 1.9
$ws.Range("P18").Value = 1.36
$ws.Range("Q18").Value = 3
$ws.Range("T18").Value = 10
$ws.Range("Z18").Value = 11
$ws.Range("AG18").Value = 23
$ws.Range("AJ18").Value = 151

# Row 22 - previously empty odds, now populated (J, K, AJ remain empty)
$ws.Range("G22").Value = 5.6
$ws.Range("H22").Value = 3.75
$ws.Range("I22").Value = 1.55
$ws.Range("L22").Value = 1.31
$ws.Range("M22").Value = 2.9
$ws.Range("N22").Value = 1.91
$ws.Range("O22").Value = 1.7
$ws.Range("P22").Value = 1.38
$ws.Range("Q22").Value = 2.57
$ws.Range("R22").Value = 1.98
$ws.Range("S22").Value = 1.65
$ws.Range("T22").Value = 13.5
$ws.Range("U22").Value = 32
$ws.Range("V22").Value = 18.5
$ws.Range("W22").Value = 120
$ws.Range("X22").Value = 65
$ws.Range("Y22").Value = 70
$ws.Range("Z22").Value = 9.25
$ws.Range("AA22").Value = 7.4
$ws.Range("AB22").Value = 19.5
$ws.Range("AC22").Value = 110
$ws.Range("AD22").Value = 5.9
$ws.Range("AE22").Value = 6.6
$ws.Range("AF22").Value = 8.25
$ws.Range("AG22").Value = 10.5
$ws.Range("AH22").Value = 13.5
$ws.Range("AI22").Value = 32

# Row 23 - previously empty odds, now fully populated
$ws.Range("G23").Value = 2.4
$ws.Range("H23").Value = 2.88
$ws.Range("I23").Value = 2.9
$ws.Range("J23").Value = 1.11
$ws.Range("K23").Value = 6.5
$ws.Range("L23").Value = 1.5
$ws.Range("M23").Value = 2.5
$ws.Range("N23").Value = 2.5
$ws.Range("O23").Value = 1.5
$ws.Range("P23").Value = 1.53
$ws.Range("Q23").Value = 2.38
$ws.Range("R23").Value = 2.1
$ws.Range("S23").Value = 1.67
$ws.Range("T23").Value = 6.5
$ws.Range("U23").Value = 11
$ws.Range("V23").Value = 11
$ws.Range("W23").Value = 26
$ws.Range("X23").Value = 23
$ws.Range("Y23").Value = 41
$ws.Range("Z23").Value = 6.5
$ws.Range("AA23").Value = 6
$ws.Range("AB23").Value = 19
$ws.Range("AC23").Value = 67
$ws.Range("AD23").Value = 7
$ws.Range("AE23").Value = 13
$ws.Range("AF23").Value = 12
$ws.Range("AG23").Value = 34
$ws.Range("AH23").Value = 29
$ws.Range("AI23").Value = 41
$ws.Range("AJ23").Value = 301
